# Update the "dSF" column (F) values for the relevant rows, per the
# "repull data, push all data, mean calculation" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -1
$ws.Range("F3").Value = 2
$ws.Range("F4").Value = -1
$ws.Range("F5").Value = 3
$ws.Range("F6").Value = -1
$ws.Range("F8").Value = -3
$ws.Range("F11").Value = 4
